# Setlist update for "latest gig" -- Mach II @ Green Man, 18 Oct 2024
# Re-keys every set's song list, re-sizes the song/total fonts, renumbers
# "Set - 4" -> "Set - 8", and drops the old "Set - 10" block entirely.

$d = $word.ActiveDocument

function Set-ParaBody($para, [string]$text, [double]$sizePt) {
    # Replace a paragraph's visible text while leaving its paragraph mark
    # alone (no trailing `r -- that would insert a whole new paragraph),
    # then size only the text run -- avoids stamping a stray <w:rPr> onto
    # <w:pPr>.
    $para.Range.Text = $text
    $full = $para.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Font.Size = $sizePt
}

# ---------------------------------------------------------------------
# 1) Global header re-brand: venue + date, every set block.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Mach II - Vogelmorn Bowling Club", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Mach II - Green Man", 2) | Out-Null
$d.Content.Find.Execute("2024-10-05", $true, $false, $false, $false, $false, `
    $true, 1, $false, "18 Oct 2024", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Drop the whole old "Set - 10" block, including the page break that
#    introduces it (paragraphs 59-65 in the original document).
# ---------------------------------------------------------------------
$pageBreakPara = $d.Paragraphs(59)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$killRange = $d.Range($pageBreakPara.Range.Start, $lastPara.Range.End)
$killRange.Delete()

# ---------------------------------------------------------------------
# 3) Old "Set - 4" (13 songs) -> new "Set - 8" (11 songs).
# ---------------------------------------------------------------------
$set4Songs = @(
    '1 : Are You Gonna Be My Girl - 3:33',
    '2 : Song 2 - 2:01',
    '3 : Fire Woman - 5:07',
    '4 : No One Knows - 4:38',
    '5 : Hash Pipe - 3:06',
    "6 : Molly's Chambers - 2:14",
    '7 : Lonely Boy - 3:13',
    '8 : Addicted To Love - 6:03',
    '9 : Brown Sugar - 3:48',
    '10 : Bliss - 4:59',
    '11 : Under the Bridge - 4:24'
)
$set4Total = 'Total Set Length: 43 min 6 sec'

$d.Paragraphs(44).Range.Text = "Set - 8"

# 13 existing song paragraphs (45-57) but only 11 needed -> delete the
# last two (56, 57) before re-keying the remaining ten.
$delRange = $d.Range($d.Paragraphs(56).Range.Start, $d.Paragraphs(57).Range.End)
$delRange.Delete()

for ($i = 0; $i -lt $set4Songs.Count; $i++) {
    Set-ParaBody $d.Paragraphs(45 + $i) $set4Songs[$i] 20
}
Set-ParaBody $d.Paragraphs(45 + $set4Songs.Count) $set4Total 16

# ---------------------------------------------------------------------
# 4) Old "Set - 3" (10 songs) -> new "Set - 3" (11 songs).
# ---------------------------------------------------------------------
$set3Songs = @(
    '1 : Tush - 2:13',
    "2 : Gimme All Your Lovin' - 4:01",
    "3 : Say It Ain't So - 4:18",
    '4 : Rebel Yell - 4:48',
    '5 : Everlong - 4:10',
    '6 : White Wedding - Pt. 1 - 4:12',
    '7 : Hard To Handle - 3:08',
    '8 : Wild Flower - 3:37',
    '9 : Are You Gonna Go My Way - 3:31',
    '10 : Wanted Dead Or Alive - 5:08',
    "11 : Sweet Child O' Mine - 5:56"
)
$set3Total = 'Total Set Length: 45 min 2 sec'

# One extra song paragraph is needed -- clone paragraph 40 (the last
# song line) so the new slot inherits matching paragraph formatting.
$d.Paragraphs(40).Range.InsertParagraphAfter()

for ($i = 0; $i -lt $set3Songs.Count; $i++) {
    Set-ParaBody $d.Paragraphs(31 + $i) $set3Songs[$i] 20
}
Set-ParaBody $d.Paragraphs(31 + $set3Songs.Count) $set3Total 16

# ---------------------------------------------------------------------
# 5) Old "Set - 2" (10 songs) -> new "Set - 2" (12 songs).
# ---------------------------------------------------------------------
$set2Songs = @(
    '1 : Crazy Little Thing Called Love - 2:43',
    '2 : Run To You - 3:53',
    '3 : Santeria - 3:02',
    '4 : Seven Nation Army - 3:52',
    '5 : Come As You Are - 3:38',
    '6 : Sweet Home Alabama - 4:43',
    '7 : White Room - 4:58',
    '8 : Honky Tonk Woman - 4:53',
    '9 : Just What I Needed - 3:45',
    '10 : Walk This Way - 3:40',
    '11 : Life in the Fast Lane - 4:46',
    '12 : Santa Monica - 3:11'
)
$set2Total = 'Total Set Length: 47 min 4 sec'

# Two extra song paragraphs are needed.
$d.Paragraphs(26).Range.InsertParagraphAfter()
$d.Paragraphs(27).Range.InsertParagraphAfter()

for ($i = 0; $i -lt $set2Songs.Count; $i++) {
    Set-ParaBody $d.Paragraphs(17 + $i) $set2Songs[$i] 20
}
Set-ParaBody $d.Paragraphs(17 + $set2Songs.Count) $set2Total 16

# ---------------------------------------------------------------------
# 6) Old "Set - 1" (10 songs) -> new "Set - 1" (11 songs).
# ---------------------------------------------------------------------
$set1Songs = @(
    '1 : April Sun in Cuba - 3:27',
    '2 : Are You Old Enough - 4:08',
    '3 : Why Does Love Do This To Me - 3:25',
    '4 : Heroes - 6:11',
    '5 : Interstate Love Song - 3:14',
    '6 : Be Mine Tonight - 6:07',
    '7 : Proud Mary -  3:07',
    '8 : Blue Lady - 3:53',
    "9 : Mary Jane's Last Dance - 4:33",
    '10 : Rain - 3:38',
    '11 : Creep - 3:58'
)
$set1Total = 'Total Set Length: 45 min 41 sec'

# One extra song paragraph is needed.
$d.Paragraphs(12).Range.InsertParagraphAfter()

for ($i = 0; $i -lt $set1Songs.Count; $i++) {
    Set-ParaBody $d.Paragraphs(3 + $i) $set1Songs[$i] 20
}
Set-ParaBody $d.Paragraphs(3 + $set1Songs.Count) $set1Total 16

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
